# Auto-generated Excel COM-interop script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain number-looking string (e.g. "4.40", "0.999").
# Excel would otherwise coerce these into real numbers (dropping trailing zeros,
# e.g. "4.40" -> 4.4), so force a Text number format first to keep them as strings,
# matching the original inline-string cell content.
$numericLookingCells = @(
    "D5", "D6", "D7", "D9", "D15", "D18", "D19", "D20", "D22", "D23", "D24", "D27", "D32", "D33", "D36", "D37", "D39", "D40", "D41", "D42", "D43", "D44", "D49", "D50", "D51"
)
foreach ($cellRef in $numericLookingCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated price / 1h-volume figures scraped for this run.
$ws.Range("D2").Value = "59.532.31"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "2.643.41"
$ws.Range("E3").Value = "  +1.43%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "537.27"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").Value = "145.34"
$ws.Range("E6").Value = "  +3.52%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +0.72%  "
$ws.Range("D9").Value = "6.66"
$ws.Range("E9").Value = "  +3.30%  "
$ws.Range("E10").Value = "  -0.41%  "
$ws.Range("E11").Value = "  +0.71%  "
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("D13").Value = "3.108.67"
$ws.Range("E13").Value = "  +1.33%  "
$ws.Range("D14").Value = "59.465.17"
$ws.Range("E14").Value = "  +0.41%  "
$ws.Range("D15").Value = "21.21"
$ws.Range("E15").Value = "  +3.27%  "
$ws.Range("D16").Value = "2.641.65"
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("E17").Value = "  +0.79%  "
$ws.Range("D18").Value = "339.93"
$ws.Range("E18").Value = "  -0.74%  "
$ws.Range("D19").Value = "4.40"
$ws.Range("E19").Value = "  +1.06%  "
$ws.Range("D20").Value = "10.40"
$ws.Range("E20").Value = "  +3.05%  "
$ws.Range("E21").Value = "  -1.24%  "
$ws.Range("D22").Value = "0.998"
$ws.Range("D23").Value = "66.99"
$ws.Range("E23").Value = "  -0.71%  "
$ws.Range("D24").Value = "0.416"
$ws.Range("E24").Value = "  +2.01%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").Value = "7.28"
$ws.Range("E27").Value = "  +0.90%  "
$ws.Range("E28").Value = "  +1.05%  "
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("E30").Value = "  +0.28%  "
$ws.Range("E31").Value = "  +0.64%  "
$ws.Range("D32").Value = "18.91"
$ws.Range("E32").Value = "  +0.55%  "
$ws.Range("D33").Value = "151.32"
$ws.Range("E33").Value = "  +1.62%  "
$ws.Range("E34").Value = "  +0.97%  "
$ws.Range("E35").Value = "  +1.71%  "
$ws.Range("D36").Value = "0.847"
$ws.Range("E36").Value = "  +2.71%  "
$ws.Range("D37").Value = "0.835"
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("E38").Value = "  -1.78%  "
$ws.Range("D39").Value = "288.27"
$ws.Range("E39").Value = "  +5.21%  "
$ws.Range("D40").Value = "3.60"
$ws.Range("E40").Value = "  +1.42%  "
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("D42").Value = "0.605"
$ws.Range("E42").Value = "  +1.45%  "
$ws.Range("D43").Value = "10.74"
$ws.Range("E43").Value = "  -0.15%  "
$ws.Range("D44").Value = "19.32"
$ws.Range("E44").Value = "  +3.89%  "
$ws.Range("E45").Value = "  +2.51%  "
$ws.Range("E46").Value = "  -1.57%  "
$ws.Range("D47").Value = "1.971.65"
$ws.Range("E47").Value = "  +1.40%  "
$ws.Range("E48").Value = "  +1.65%  "
$ws.Range("D49").Value = "4.57"
$ws.Range("E49").Value = "  +1.66%  "
$ws.Range("D50").Value = "18.30"
$ws.Range("E50").Value = "  +0.29%  "
$ws.Range("D51").Value = "110.87"
$ws.Range("E51").Value = "  +0.21%  "
